# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" on the Overview,
#   zh-cn and de-de sheets (the cells that hold the current hand-off status).
# - The corresponding status columns narrow to fit the new (shorter) text,
#   matching Excel's "best fit" column width after the content change.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns("E").ColumnWidth = 12.5
$wsOverview.Columns("F").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns("C").ColumnWidth = 12.5

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns("C").ColumnWidth = 12.5
